$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1898016997167139
$ws.Range("C2").Value = 0.5864022662889519
$ws.Range("J2").Value = 0.0169971671388102
$ws.Range("P2").Value = 0.1359773371104816
$ws.Range("S2").Value = 0.07082152974504249
$ws.Range("C3").Value = 0.04545454545454546
$ws.Range("J3").Value = 0.004545454545454545
$ws.Range("P3").Value = 0.7181818181818181
$ws.Range("S3").Value = 0.2318181818181818
$ws.Range("J4").Value = 0.01886792452830189
$ws.Range("P4").Value = 0.6792452830188679
$ws.Range("S4").Value = 0.3018867924528302
$ws.Range("B6").Value = 0.0546875
$ws.Range("D6").Value = 0.01171875
$ws.Range("F6").Value = 0.05859375
$ws.Range("J6").Value = 0.234375
$ws.Range("O6").Value = 0.03125
$ws.Range("Q6").Value = 0.1953125
$ws.Range("R6").Value = 0.046875
$ws.Range("S6").Value = 0.3671875
$ws.Range("B7").Value = 0.0815450643776824
$ws.Range("D7").Value = 0.01716738197424893
$ws.Range("F7").Value = 0.05579399141630902
$ws.Range("J7").Value = 0.1459227467811159
$ws.Range("O7").Value = 0.03004291845493562
$ws.Range("Q7").Value = 0.1545064377682404
$ws.Range("R7").Value = 0.06008583690987124
$ws.Range("S7").Value = 0.4549356223175965
$ws.Range("B8").Value = 0.1351931330472103
$ws.Range("D8").Value = 0.01502145922746781
$ws.Range("E8").Value = 0.002145922746781116
$ws.Range("F8").Value = 0.04291845493562232
$ws.Range("J8").Value = 0.09442060085836911
$ws.Range("O8").Value = 0.01716738197424893
$ws.Range("Q8").Value = 0.1974248927038627
$ws.Range("R8").Value = 0.06437768240343347
$ws.Range("S8").Value = 0.4313304721030043
$ws.Range("B9").Value = 0.1101321585903084
$ws.Range("D9").Value = 0.02202643171806168
$ws.Range("F9").Value = 0.06607929515418502
$ws.Range("J9").Value = 0.1277533039647577
$ws.Range("O9").Value = 0.03083700440528634
$ws.Range("Q9").Value = 0.1497797356828194
$ws.Range("R9").Value = 0.05726872246696035
$ws.Range("S9").Value = 0.4361233480176211
$ws.Range("B10").Value = 0.1268482490272373
$ws.Range("D10").Value = 0.02801556420233463
$ws.Range("E10").Value = 0.0007782101167315176
$ws.Range("F10").Value = 0.08560311284046693
$ws.Range("J10").Value = 0.09961089494163425
$ws.Range("O10").Value = 0.02101167315175097
$ws.Range("Q10").Value = 0.1937743190661479
$ws.Range("R10").Value = 0.0536964980544747
$ws.Range("S10").Value = 0.3906614785992218
$ws.Range("G11").Value = 0.152112676056338
$ws.Range("J11").Value = 0.08732394366197183
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.5464788732394367
$ws.Range("S11").Value = 0.01408450704225352
$ws.Range("G12").Value = 0.7537688442211056
$ws.Range("J12").Value = 0.2010050251256282
$ws.Range("K12").Value = 0.005025125628140704
$ws.Range("L12").Value = 0.01507537688442211
$ws.Range("S12").Value = 0.02512562814070352
$ws.Range("G13").Value = 0.7954545454545454
$ws.Range("J13").Value = 0.2045454545454546
$ws.Range("F15").Value = 0.02553191489361702
$ws.Range("H15").Value = 0.1191489361702128
$ws.Range("I15").Value = 0.08085106382978724
$ws.Range("J15").Value = 0.3276595744680851
$ws.Range("K15").Value = 0.09361702127659574
$ws.Range("M15").Value = 0.01276595744680851
$ws.Range("O15").Value = 0.0425531914893617
$ws.Range("S15").Value = 0.2978723404255319
$ws.Range("F16").Value = 0.01271186440677966
$ws.Range("H16").Value = 0.1440677966101695
$ws.Range("I16").Value = 0.1101694915254237
$ws.Range("J16").Value = 0.4067796610169492
$ws.Range("K16").Value = 0.1101694915254237
$ws.Range("M16").Value = 0.01271186440677966
$ws.Range("O16").Value = 0.04661016949152542
$ws.Range("S16").Value = 0.1567796610169492
$ws.Range("F17").Value = 0.02178649237472767
$ws.Range("H17").Value = 0.1895424836601307
$ws.Range("I17").Value = 0.1089324618736384
$ws.Range("J17").Value = 0.4183006535947713
$ws.Range("K17").Value = 0.08496732026143791
$ws.Range("M17").Value = 0.01742919389978214
$ws.Range("O17").Value = 0.04139433551198257
$ws.Range("S17").Value = 0.1176470588235294
$ws.Range("F18").Value = 0.03597122302158273
$ws.Range("H18").Value = 0.223021582733813
$ws.Range("I18").Value = 0.1079136690647482
$ws.Range("J18").Value = 0.3597122302158273
$ws.Range("K18").Value = 0.1079136690647482
$ws.Range("M18").Value = 0.02158273381294964
$ws.Range("O18").Value = 0.06474820143884892
$ws.Range("S18").Value = 0.07913669064748201
$ws.Range("F19").Value = 0.01997146932952924
$ws.Range("H19").Value = 0.2097004279600571
$ws.Range("I19").Value = 0.08416547788873038
$ws.Range("J19").Value = 0.3616262482168331
$ws.Range("K19").Value = 0.1283880171184023
$ws.Range("M19").Value = 0.02068473609129814
$ws.Range("O19").Value = 0.07489300998573467
$ws.Range("S19").Value = 0.1005706134094151
